# Apply the 2022-06-03 Fonds de solidarite data refresh.
# For a set of rows, update column C (nombre_aides) and column E (montant_total)
# to their new values, leaving column D (nombre_entreprises) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 10;  C = 278199; E = 1752276968 },
    @{ Row = 19;  C = 108920; E = 344665282 },
    @{ Row = 115; C = 17556;  E = 38619107 },
    @{ Row = 134; C = 5678;   E = 17162226 },
    @{ Row = 152; C = 126047; E = 715970384 },
    @{ Row = 168; C = 285023; E = 1210861599 },
    @{ Row = 169; C = 562613; E = 1285058739 },
    @{ Row = 170; C = 367415; E = 2846304604 },
    @{ Row = 171; C = 115169; E = 447140114 },
    @{ Row = 173; C = 54392;  E = 151908108 },
    @{ Row = 174; C = 357255; E = 1018552234 },
    @{ Row = 175; C = 125559; E = 813138243 },
    @{ Row = 177; C = 96761;  E = 174754897 },
    @{ Row = 179; C = 235723; E = 812710859 },
    @{ Row = 267; C = 84977;  E = 156521773 },
    @{ Row = 313; C = 220650; E = 1371052430 },
    @{ Row = 317; C = 103584; E = 303300964 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 3).Value = $u.C
    $ws.Cells.Item($r, 5).Value = $u.E
}

$wb.Save()
